# Daily attendance processing - 2026-01-14 23:56:19
# Swap the order of "System" and the recorder's email address inside the
# "Recorded By" column (G) so that entries read "<email>, System"
# instead of "System, <email>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colIndex = 7  # Column G = "Recorded By"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colIndex)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
